$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list with latest values (text-formatted cells to preserve original string typing)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.661.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.728.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.28"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.63"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.51%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.212.54"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.67"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.543.02"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.733.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.72"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.46"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.520"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.30"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.33"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0900"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.13"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.04%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +11.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.87"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.86"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.59%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.968"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "344.20"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.09%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.07"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.72"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.99"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0581"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.625"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0997"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.07"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.93"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.82%  "
